# Generate Report for Handback
#
# The df08d9e8-... localization file has now been handed back (in sync
# with en-US). Update its status from "Ready for handoff" to
# "Handed back: in sync with en-US" on the Overview sheet as well as the
# per-locale (zh-cn / de-de) detail sheets, and record the new "Latest
# Handback DateTime" timestamps on the detail sheets.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: row 3 is the df08d9e8-...md file ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Handed back: in sync with en-US"
$overview.Range("C3").Value = "Handed back: in sync with en-US"

# --- zh-cn detail sheet: row 3 is the df08d9e8-...md file ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Handed back: in sync with en-US"
$zhcn.Range("H3").Value = "2016-03-25 03:00:58"

# --- de-de detail sheet: row 3 is the df08d9e8-...md file ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Handed back: in sync with en-US"
$dede.Range("H3").Value = "2016-03-25 03:01:06"
